$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MirroMe Voorbeeldredenering")

# Text edits - order matters for shared-string table append order.
$ws.Range("E3").Value = "Voor de toepassing van het bij of krachtens de WOR bepaalde wordt [de bestuurder] geacht niet te behoren tot de in [de onderneming] werkzame personen."
$ws.Range("E4").Value = "Verkiesbaar tot lid van de ondernemingsraad zijn de personen die gedurende ten minste een jaar in [de onderneming] werkzaam zijn geweest."
$ws.Range("E5").Value = "[de bestuurder] is niet verkiesbaar tot lid van de ondernemingsraad van [de onderneming]."
$ws.Range("E7").Value = "[de bestuurder] kan geen voorzitter zijn van de ondernemingsraad van [de onderneming]."
$ws.Range("E12").Value = "[de onderneming] is een onderneming in de zin van de Wet op de ondernemingsraden, omdat het een in de maatschappij als zelfstandige eenheid optredend organisatorisch verband is, waarin krachtens arbeidsovereenkomst of krachtens publiekrechtelijke aanstelling arbeid wordt verricht."
$ws.Range("C13").Value = "de bestuurder"
$ws.Range("C12").Value = "de onderneming"
$ws.Range("D12").Value = "Albert Heijn"
$ws.Range("D13").Value = "Piet van der Kluns"
$ws.Range("E13").Value = "[de bestuurder] is bestuurder in de zin van de Wet op de ondernemingsraden, omdat hij in [de onderneming] de hoogste zeggenschap uitoefent bij de leiding van de arbeid."

# Selection moves to A2 on the first sheet.
$ws.Activate()
$ws.Range("A2").Select()
